$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.997.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.583.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.59%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.543'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.608.63'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.115'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.72%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.361'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000184'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.059.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.969.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.611.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '66.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '640.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000106'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +11.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.689.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.90'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.139'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.84%  '
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.98'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.34%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '156.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.91%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '19.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.374'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.83'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '162.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.14%  '
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.636'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.77%  '
